$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.848.91'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '2.321.99'
$ws.Range("E3").Value = '  +4.03%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = "'97.94"
$ws.Range("E5").Value = '  +5.67%  '
$ws.Range("D6").Value = "'270.84"
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").Value = "'45.97"
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").Value = "'0.0951"
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("D12").Value = "'8.12"
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("D14").Value = '2.673.29'
$ws.Range("E14").Value = '  +4.08%  '
$ws.Range("D15").Value = "'15.51"
$ws.Range("E15").Value = '  +2.86%  '
$ws.Range("E16").Value = '  +8.00%  '
$ws.Range("D17").Value = '2.332.16'
$ws.Range("E17").Value = '  +4.17%  '
$ws.Range("D18").Value = '43.824.73'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E19").Value = '  +5.25%  '
$ws.Range("E20").Value = '  +6.81%  '
$ws.Range("D21").Value = "'72.83"
$ws.Range("E21").Value = '  +3.42%  '
$ws.Range("D22").Value = "'239.47"
$ws.Range("E22").Value = '  +2.65%  '
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("D24").Value = "'9.45"
$ws.Range("E24").Value = '  +5.13%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = "'2.52"
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").Value = "'11.34"
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = "'3.47"
$ws.Range("E28").Value = '  -2.12%  '
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").Value = "'38.35"
$ws.Range("E30").Value = '  -4.68%  '
$ws.Range("D31").Value = "'22.38"
$ws.Range("E31").Value = '  +7.49%  '
$ws.Range("D32").Value = "'175.18"
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("D33").Value = "'0.0908"
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  +3.13%  '
$ws.Range("D36").Value = "'0.0360"
$ws.Range("E36").Value = '  +2.52%  '
$ws.Range("E37").Value = '  -2.58%  '
$ws.Range("D38").Value = "'4.41"
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("E39").Value = '  -5.95%  '
$ws.Range("E40").Value = '  +10.44%  '
$ws.Range("D41").Value = "'2.35"
$ws.Range("E41").Value = '  +8.07%  '
$ws.Range("E42").Value = '  +19.28%  '
$ws.Range("D43").Value = "'12.22"
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("D44").Value = "'9.21"
$ws.Range("E44").Value = '  +9.85%  '
$ws.Range("D45").Value = "'62.18"
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("D46").Value = "'5.39"
$ws.Range("E46").Value = '  +1.02%  '
$ws.Range("E47").Value = '  +4.23%  '
$ws.Range("D48").Value = "'100.36"
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("E50").Value = '  +16.68%  '
$ws.Range("D51").Value = '2.550.49'
$ws.Range("E51").Value = '  +3.82%  '
